# FSD-Automation.docx template update
#   - wrap a few template variables with the Jinja "|escape" filter
#   - drop the stray "_GoBack" bookmark
#   - nudge the version-history table's cached column grid
#   - same "|escape" wrap for the header's {{group}} variable

$d = $word.ActiveDocument

# --- word/document.xml ---------------------------------------------------

# {{service}} -> {{service|escape}}  (Service / PIC table, row 2)
$d.Content.Find.Execute(
    "{{service}}", $false, $false, $false, $false, $false,
    $true, 1, $false, "{{service|escape}}", 2) | Out-Null

# {{activity}} -> {{activity|escape}}  (version history table header row)
$d.Content.Find.Execute(
    "{{activity}}", $false, $false, $false, $false, $false,
    $true, 1, $false, "{{activity|escape}}", 2) | Out-Null

# {{scope}} -> {{scope|escape}}  (1.2 Scope section)
$d.Content.Find.Execute(
    "{{scope}}", $false, $false, $false, $false, $false,
    $true, 1, $false, "{{scope|escape}}", 2) | Out-Null

# Drop the leftover "_GoBack" bookmark (Word keeps re-creating/removing this
# one between edit sessions; it carries no content).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-balance the cached grid of the Version/Date/Author/Activity table —
# total width is unchanged (9030 dxa), only the inner column boundaries move.
$historyTable = $d.Tables(2)
$historyTable.Columns(2).Width = 1448 / 20.0
$historyTable.Columns(3).Width = 2936 / 20.0
$historyTable.Columns(4).Width = 3641 / 20.0

# --- word/header1.xml ------------------------------------------------------

foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers(1)
    $hdr.Range.Find.Execute(
        "{{group}}", $false, $false, $false, $false, $false,
        $true, 1, $false, "{{group|escape}}", 2) | Out-Null
}
